$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.683.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.62%  "

$ws.Range("D3").Value = "'1.605.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'212.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("D6").Value = "'0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'27.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.54%  "

$ws.Range("D9").Value = "'43.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("D10").Value = "'0.251"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.27%  "

$ws.Range("D11").Value = "'0.0600"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.17%  "

$ws.Range("D12").Value = "'0.0908"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("D13").Value = "'1.837.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.83%  "

$ws.Range("D14").Value = "'1.611.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("D15").Value = "'29.680.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.51%  "

$ws.Range("D16").Value = "'0.537"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.89%  "

$ws.Range("D17").Value = "'3.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("D18").Value = "'63.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.90%  "

$ws.Range("D19").Value = "'241.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.93%  "

$ws.Range("D20").Value = "'7.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.94%  "

$ws.Range("D21").Value = "'0.0₃0693"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.60%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Value = "'3.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("D24").Value = "'9.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "

$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("D26").Value = "'155.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.06%  "

$ws.Range("E27").Value = "  +3.95%  "

$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("E29").Value = "  +2.66%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").Value = "  +3.83%  "

$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("D34").Value = "'1.433.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.22%  "

$ws.Range("E35").Value = "  +4.19%  "

$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("E37").Value = "  +4.84%  "

$ws.Range("D38").Value = "'2.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.04%  "

$ws.Range("D39").Value = "'2.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").Value = "'0.0165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "

$ws.Range("D41").Value = "'0.539"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.05%  "

$ws.Range("D42").Value = "'1.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.85%  "

$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'54.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +27.39%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0488"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.25%  "

$ws.Range("D45").Value = "'0.801"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.48%  "

$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("D47").Value = "'65.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.94%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'5.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'0.942"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.89%  "

$ws.Range("D50").Value = "'1.747.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.98%  "

$ws.Range("D51").Value = "'86.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
